$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'46.468.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = "'2.598.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +9.98%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'305.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.30%  '
$ws.Range("D6").Value = "'100.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("D7").Value = "'0.601"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.42%  '
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +13.18%  '
$ws.Range("D10").Value = "'38.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +11.92%  '
$ws.Range("D11").Value = "'0.0838"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.65%  '
$ws.Range("D12").Value = "'8.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +14.76%  '
$ws.Range("D13").Value = "'2.993.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.93%  '
$ws.Range("E14").Value = '  +1.97%  '
$ws.Range("D15").Value = "'2.595.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +10.18%  '
$ws.Range("D16").Value = "'0.903"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +10.65%  '
$ws.Range("E17").Value = '  +9.19%  '
$ws.Range("D18").Value = "'46.656.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.31%  '
$ws.Range("D19").Value = "'13.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.75%  '
$ws.Range("E20").Value = '  +4.34%  '
$ws.Range("E21").Value = '  +9.76%  '
$ws.Range("D22").Value = "'71.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.40%  '
$ws.Range("D23").Value = "'259.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.36%  '
$ws.Range("D24").Value = "'2.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.82%  '
$ws.Range("D25").Value = "'2.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.11%  '
$ws.Range("D26").Value = "'28.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +33.55%  '
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").Value = "'10.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.14%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").Value = "'39.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.96%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = "'2.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.68%  '
$ws.Range("D31").Value = "'3.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.64%  '
$ws.Range("D32").Value = "'6.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.17%  '
$ws.Range("E33").Value = '  +21.87%  '
$ws.Range("E34").Value = '  +5.19%  '
$ws.Range("E35").Value = '  +7.42%  '
$ws.Range("D36").Value = "'149.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.38%  '
$ws.Range("D37").Value = "'0.117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.91%  '
$ws.Range("E39").Value = '  +5.90%  '
$ws.Range("D40").Value = "'15.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.66%  '
$ws.Range("D41").Value = "'3.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +12.72%  '
$ws.Range("D43").Value = "'2.036.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.61%  '
$ws.Range("D44").Value = "'18.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +29.15%  '
$ws.Range("D45").Value = "'0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").Value = "'91.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").Value = "'1.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.83%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Value = "'9.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.54%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = "'109.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +11.22%  '
$ws.Range("D50").Value = "'0.202"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.70%  '
$ws.Range("D51").Value = "'2.851.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.91%  '
